$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 252-255 continue the daily series with the same style as the
# preceding row (A251) - style index 2 (date format, centered, bordered, bold).
$ws.Range("A251").Copy()
$ws.Range("A252:A255").PasteSpecial(-4122)

$ws.Cells.Item(252, 1).Value = 44326
$ws.Cells.Item(252, 2).Value = 0
$ws.Cells.Item(252, 3).Value = 0
$ws.Cells.Item(252, 4).Value = 0

$ws.Cells.Item(253, 1).Value = 44327
$ws.Cells.Item(253, 2).Value = 0
$ws.Cells.Item(253, 3).Value = 0
$ws.Cells.Item(253, 4).Value = 0

$ws.Cells.Item(254, 1).Value = 44328
$ws.Cells.Item(254, 2).Value = 1
$ws.Cells.Item(254, 3).Value = 1
$ws.Cells.Item(254, 4).Value = 26.76659528907923

$ws.Cells.Item(255, 1).Value = 44329
$ws.Cells.Item(255, 2).Value = 1
$ws.Cells.Item(255, 3).Value = 2
$ws.Cells.Item(255, 4).Value = 53.53319057815846
